$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3056.9614
$ws.Range("I64").Value = 2965.111
$ws.Range("J64").Value = 3263.625
$ws.Range("K64").Value = 2965.111
$ws.Range("L64").Value = 3263.625
$ws.Range("M64").Value = -2717.111
$ws.Range("N64").Value = -3759.625
$ws.Range("H67").Value = 3056.9614
$ws.Range("I67").Value = 2965.111
$ws.Range("J67").Value = 3263.625
$ws.Range("K67").Value = 2965.111
$ws.Range("L67").Value = 3263.625
$ws.Range("M67").Value = -2107.111
$ws.Range("N67").Value = -4979.625
$ws.Range("H74").Value = 4300.4
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 4222.6665
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 4222.6665
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -6094.6665
$ws.Range("H77").Value = 4300.4
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 4222.6665
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 21113.3325
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -30473.3325
$ws.Range("H106").Value = 926
$ws.Range("I106").Value = 926
$ws.Range("K106").Value = 926
$ws.Range("M106").Value = -295

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 5000
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").Value = 25000
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 2458.4866
$ws.Range("I74").Value = 2443.5806
$ws.Range("K74").Value = 2443.5806
$ws.Range("M74").Value = -1569.5806
$ws.Range("H77").Value = 2458.4866
$ws.Range("I77").Value = 2443.5806
$ws.Range("K77").Value = 12217.903
$ws.Range("M77").Value = -7849.902999999998
$ws.Range("H102").Value = 942
$ws.Range("I102").Value = 942
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 942
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 680
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0
$ws.Range("H132").Value = 2217.1147
$ws.Range("I132").Value = 1978.3846
$ws.Range("J132").Value = 3596.4443
$ws.Range("K132").Value = 5935.1538
$ws.Range("L132").Value = 10789.3329
$ws.Range("M132").Value = -3405.1538
$ws.Range("N132").Value = -15849.3329

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2115.5881
$ws.Range("I105").Value = 1687
$ws.Range("K105").Value = 1687
$ws.Range("M105").Value = 60
$ws.Range("H107").Value = 1481.8572
$ws.Range("I107").Value = 1181.174
$ws.Range("J107").Value = 2058.1667
$ws.Range("K107").Value = 1181.174
$ws.Range("L107").Value = 2058.1667
$ws.Range("M107").Value = 738.826
$ws.Range("N107").Value = -5898.1667
$ws.Range("H117").Value = 42000
$ws.Range("J117").Value = 42000
$ws.Range("L117").Value = 42000
$ws.Range("N117").Value = -51178

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3174.7144
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 3438.4443
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 3438.4443
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -4686.4443
$ws.Range("H65").Value = 3174.7144
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 3438.4443
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 17192.2215
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -23432.2215
$ws.Range("H107").Value = 831.2692
$ws.Range("I107").Value = 701.5238000000001
$ws.Range("J107").Value = 1376.2
$ws.Range("K107").Value = 701.5238000000001
$ws.Range("L107").Value = 1376.2
$ws.Range("M107").Value = 1218.4762
$ws.Range("N107").Value = -5216.2
$ws.Range("H132").Value = 2799.75
$ws.Range("I132").Value = 2157
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 6471
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -3941
$ws.Range("N132").Value = -16158.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 4800
$ws.Range("I46").Value = 4800
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4800
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4644
$ws.Range("H70").Value = 103679730
$ws.Range("I70").Value = 207355200
$ws.Range("J70").Value = 4262.5
$ws.Range("K70").Value = 207355200
$ws.Range("L70").Value = 4262.5
$ws.Range("M70").Value = -207354930
$ws.Range("N70").Value = -4802.5
$ws.Range("H73").Value = 103679730
$ws.Range("I73").Value = 207355200
$ws.Range("J73").Value = 4262.5
$ws.Range("K73").Value = 207355200
$ws.Range("L73").Value = 4262.5
$ws.Range("M73").Value = -207354264
$ws.Range("N73").Value = -6134.5
$ws.Range("H80").Value = 2465.5
$ws.Range("I80").Value = 2625
$ws.Range("J80").Value = 2385.75
$ws.Range("K80").Value = 2625
$ws.Range("L80").Value = 2385.75
$ws.Range("M80").Value = -1627
$ws.Range("N80").Value = -4381.75
$ws.Range("H83").Value = 2465.5
$ws.Range("I83").Value = 2625
$ws.Range("J83").Value = 2385.75
$ws.Range("K83").Value = 13125
$ws.Range("L83").Value = 11928.75
$ws.Range("M83").Value = -8133
$ws.Range("N83").Value = -21912.75
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
$ws.Range("H132").Value = 3801.9429
$ws.Range("I132").Value = 3697.625
$ws.Range("J132").Value = 4029.5454
$ws.Range("K132").Value = 11092.875
$ws.Range("L132").Value = 12088.6362
$ws.Range("M132").Value = -8562.875
$ws.Range("N132").Value = -17148.6362

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7468269
$ws.Range("I132").Value = 10210490
$ws.Range("J132").Value = 3334.889
$ws.Range("K132").Value = 30631470
$ws.Range("L132").Value = 10004.667
$ws.Range("M132").Value = -30628940
$ws.Range("N132").Value = -15064.667
